# Listas sem duplicação de professores
# Replace cells that still contained duplicate/bracketed teacher lists with "-"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cells = @("B2", "D2", "B3", "D3", "E3", "F3", "B4", "D4", "E4", "F4", "B6", "D6", "F6", "B7", "C7", "D7", "F7", "B8", "D8")

foreach ($cell in $cells) {
    $ws.Range($cell).Value = "-"
}
